# Add a new week (row 21) of comprehension-scores data to the log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week number
$ws.Range("A21").Value = 20

# Time spent that week (stored as an Excel time fraction of a day,
# formatted like the other "under a day" rows such as B9/B13/B15).
$ws.Range("B21").Value = 0.6492013888888889
$ws.Range("B21").NumberFormat = "h:mm:ss"

# Running total formula, matching the pattern used by the rows above it.
$ws.Range("C21").Formula = "=SUM(B2:B21)+1.2708333333"

# What I watched/read that week.
$ws.Range("D21").Value = "[10 Cosas Respondidas por Mojang! #7](https://youtu.be/jl55cLP5IAw) (Audiovisual, Spanish, New):45;"

# Move the active selection to C22, like the source workbook shows.
[void]$ws.Range("C22").Select()
